$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns C (td_sim_1) and D (record_atd), rows 2-21,
# plus the summary cell C22 (average of td_sim_1).

$updates = @(
    @{ Row = 2;  C = 82;                D = 83 },
    @{ Row = 3;  C = 84;                D = 85.5 },
    @{ Row = 4;  C = 12;                D = 10 },
    @{ Row = 5;  C = 228;               D = 221.5 },
    @{ Row = 6;  C = 19;                D = 12 },
    @{ Row = 7;  C = 290;               D = 329 },
    @{ Row = 8;  C = 21;                D = 12 },
    @{ Row = 9;  C = 30;                D = 30 },
    @{ Row = 10; C = 23;                D = 16.5 },
    @{ Row = 11; C = 6;                 D = 9 },
    @{ Row = 12; C = 54;                D = 59.5 },
    @{ Row = 13; C = 229;               D = 223 },
    @{ Row = 14; C = 72;                D = 76 },
    @{ Row = 15; C = 33;                D = 28 },
    @{ Row = 16; C = 38;                D = 32 },
    @{ Row = 17; C = 57;                D = 61.5 },
    @{ Row = 18; C = 56;                D = 64.5 },
    @{ Row = 19; C = 7;                 D = 10 },
    @{ Row = 20; C = 155;               D = 158 },
    @{ Row = 21; C = 52;                D = 57 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

# Summary row: average td_sim_1 value
$ws.Cells.Item(22, 3).Value = 77.40000000000001
